$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATMI ligand-receptor pair metrics (Zp3-Mertk) with refreshed TPM values.
$ws.Range("G2").Value = 0.08261866666666666
$ws.Range("H2").Value = 0.247856
$ws.Range("I2").Value = 0.5806915650061265
$ws.Range("J2").Value = 0.5806915650061265
$ws.Range("M2").Value = 11.01658666666667
$ws.Range("N2").Value = 33.04976
$ws.Range("O2").Value = 0.1837409300120545
$ws.Range("P2").Value = 0.1837409300120545
$ws.Range("Q2").Value = 0.9101757016177777
$ws.Range("R2").Value = 8.191581314559999
$ws.Range("S2").Value = 0.1066968082043811
$ws.Range("T2").Value = 0.1066968082043811
$ws.Range("G3").Value = 0.08261866666666666
$ws.Range("H3").Value = 0.247856
$ws.Range("I3").Value = 0.5806915650061265
$ws.Range("J3").Value = 0.5806915650061265
$ws.Range("O3").Value = 0.06509859443116503
$ws.Range("P3").Value = 0.06509859443116503
$ws.Range("Q3").Value = 0.3224712036497777
$ws.Range("R3").Value = 2.902240832848
$ws.Range("S3").Value = 0.03780220467993233
$ws.Range("T3").Value = 0.03780220467993233
$ws.Range("G4").Value = 0.08261866666666666
$ws.Range("H4").Value = 0.247856
$ws.Range("I4").Value = 0.5806915650061265
$ws.Range("J4").Value = 0.5806915650061265
$ws.Range("M4").Value = 4.001997666666667
$ws.Range("N4").Value = 12.005993
$ws.Range("O4").Value = 0.06674760480978428
$ws.Range("P4").Value = 0.06674760480978428
$ws.Range("Q4").Value = 0.3306397112231111
$ws.Range("R4").Value = 2.975757401008
$ws.Range("S4").Value = 0.03875977109740409
$ws.Range("T4").Value = 0.03875977109740409
$ws.Range("G5").Value = 0.08261866666666666
$ws.Range("H5").Value = 0.247856
$ws.Range("I5").Value = 0.5806915650061265
$ws.Range("J5").Value = 0.5806915650061265
$ws.Range("M5").Value = 41.03546066666667
$ws.Range("N5").Value = 123.106382
$ws.Range("O5").Value = 0.6844128707469962
$ws.Range("P5").Value = 0.6844128707469963
$ws.Range("Q5").Value = 3.390295046332444
$ws.Range("R5").Value = 30.512655416992
$ws.Range("S5").Value = 0.397432781024409
$ws.Range("T5").Value = 0.397432781024409
$ws.Range("I6").Value = 0.4193084349938734
$ws.Range("J6").Value = 0.4193084349938734
$ws.Range("M6").Value = 11.01658666666667
$ws.Range("N6").Value = 33.04976
$ws.Range("O6").Value = 0.1837409300120545
$ws.Range("P6").Value = 0.1837409300120545
$ws.Range("Q6").Value = 0.6572238551644445
$ws.Range("R6").Value = 5.915014696479999
$ws.Range("S6").Value = 0.07704412180767339
$ws.Range("T6").Value = 0.07704412180767341
$ws.Range("I7").Value = 0.4193084349938734
$ws.Range("J7").Value = 0.4193084349938734
$ws.Range("O7").Value = 0.06509859443116503
$ws.Range("P7").Value = 0.06509859443116503
$ws.Range("S7").Value = 0.02729638975123269
$ws.Range("T7").Value = 0.02729638975123269
$ws.Range("I8").Value = 0.4193084349938734
$ws.Range("J8").Value = 0.4193084349938734
$ws.Range("M8").Value = 4.001997666666667
$ws.Range("N8").Value = 12.005993
$ws.Range("O8").Value = 0.06674760480978428
$ws.Range("P8").Value = 0.06674760480978428
$ws.Range("Q8").Value = 0.2387498427987778
$ws.Range("R8").Value = 2.148748585189
$ws.Range("S8").Value = 0.02798783371238018
$ws.Range("T8").Value = 0.02798783371238018
$ws.Range("I9").Value = 0.4193084349938734
$ws.Range("J9").Value = 0.4193084349938734
$ws.Range("M9").Value = 41.03546066666667
$ws.Range("N9").Value = 123.106382
$ws.Range("O9").Value = 0.6844128707469962
$ws.Range("P9").Value = 0.6844128707469963
$ws.Range("Q9").Value = 2.448079833965111
$ws.Range("R9").Value = 22.032718505686
$ws.Range("S9").Value = 0.2869800897225871
$ws.Range("T9").Value = 0.2869800897225872
